$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Lin7c"
$ws.Range("C2").Value = "Htr2c"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 8.190874333333333
$ws.Range("H2").Value = 24.572623
$ws.Range("I2").Value = 0.1295676603100907
$ws.Range("J2").Value = 0.1394023768750731
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.0002613333333333333
$ws.Range("N2").Value = 0.000784
$ws.Range("O2").Value = 0.2316100443131462
$ws.Range("P2").Value = 0.2316100443131462
$ws.Range("Q2").Value = 0.002140548492444444
$ws.Range("R2").Value = 0.019264936432
$ws.Range("S2").Value = 0.03000917154597078
$ws.Range("T2").Value = 0.0322869906853936

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Lin7c"
$ws.Range("C3").Value = "Htr2c"
$ws.Range("D3").Value = "Inflammatory-Mac"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 8.190874333333333
$ws.Range("H3").Value = 24.572623
$ws.Range("I3").Value = 0.1295676603100907
$ws.Range("J3").Value = 0.1394023768750731
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.000867
$ws.Range("N3").Value = 0.002601
$ws.Range("O3").Value = 0.7683899556868538
$ws.Range("P3").Value = 0.7683899556868538
$ws.Range("Q3").Value = 0.007101488047000001
$ws.Range("R3").Value = 0.063913392423
$ws.Range("S3").Value = 0.09955848876411992
$ws.Range("T3").Value = 0.1071153861896795

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Lin7c"
$ws.Range("C4").Value = "Htr2c"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 8.554364999999999
$ws.Range("H4").Value = 25.663095
$ws.Range("I4").Value = 0.1353175513849534
$ws.Range("J4").Value = 0.1455887082535228
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.0002613333333333333
$ws.Range("N4").Value = 0.000784
$ws.Range("O4").Value = 0.2316100443131462
$ws.Range("P4").Value = 0.2316100443131462
$ws.Range("Q4").Value = 0.00223554072
$ws.Range("R4").Value = 0.02011986648
$ws.Range("S4").Value = 0.03134090407261549
$ws.Range("T4").Value = 0.03371980717009214

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Lin7c"
$ws.Range("C5").Value = "Htr2c"
$ws.Range("D5").Value = "Inflammatory-Mac"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 8.554364999999999
$ws.Range("H5").Value = 25.663095
$ws.Range("I5").Value = 0.1353175513849534
$ws.Range("J5").Value = 0.1455887082535228
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.000867
$ws.Range("N5").Value = 0.002601
$ws.Range("O5").Value = 0.7683899556868538
$ws.Range("P5").Value = 0.7683899556868538
$ws.Range("Q5").Value = 0.007416634455
$ws.Range("R5").Value = 0.066749710095
$ws.Range("S5").Value = 0.1039766473123379
$ws.Range("T5").Value = 0.1118689010834307

# Row 6
$ws.Range("A6").Value = "Inflammatory-Mac"
$ws.Range("B6").Value = "Lin7c"
$ws.Range("C6").Value = "Htr2c"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 15.38058366666667
$ws.Range("H6").Value = 46.141751
$ws.Range("I6").Value = 0.2432983536059943
$ws.Range("J6").Value = 0.2617656960177911
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.0002613333333333333
$ws.Range("N6").Value = 0.000784
$ws.Range("O6").Value = 0.2316100443131462
$ws.Range("P6").Value = 0.2316100443131462
$ws.Range("Q6").Value = 0.004019459198222222
$ws.Range("R6").Value = 0.036175132784
$ws.Range("S6").Value = 0.05635034245999985
$ws.Range("T6").Value = 0.06062756445434217

# Row 7
$ws.Range("A7").Value = "Inflammatory-Mac"
$ws.Range("B7").Value = "Lin7c"
$ws.Range("C7").Value = "Htr2c"
$ws.Range("D7").Value = "Inflammatory-Mac"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 15.38058366666667
$ws.Range("H7").Value = 46.141751
$ws.Range("I7").Value = 0.2432983536059943
$ws.Range("J7").Value = 0.2617656960177911
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.000867
$ws.Range("N7").Value = 0.002601
$ws.Range("O7").Value = 0.7683899556868538
$ws.Range("P7").Value = 0.7683899556868538
$ws.Range("Q7").Value = 0.013334966039
$ws.Range("R7").Value = 0.120014694351
$ws.Range("S7").Value = 0.1869480111459944
$ws.Range("T7").Value = 0.201138131563449

# Row 8
$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Lin7c"
$ws.Range("C8").Value = "Htr2c"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 13.379706
$ws.Range("H8").Value = 26.759412
$ws.Range("I8").Value = 0.2116473933682475
$ws.Range("J8").Value = 0.1518081987657302
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.0002613333333333333
$ws.Range("N8").Value = 0.000784
$ws.Range("O8").Value = 0.2316100443131462
$ws.Range("P8").Value = 0.2316100443131462
$ws.Range("Q8").Value = 0.003496563167999999
$ws.Range("R8").Value = 0.020979379008
$ws.Range("S8").Value = 0.0490196621567817
$ws.Range("T8").Value = 0.03516030364322969

# Row 9
$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Lin7c"
$ws.Range("C9").Value = "Htr2c"
$ws.Range("D9").Value = "Inflammatory-Mac"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 13.379706
$ws.Range("H9").Value = 26.759412
$ws.Range("I9").Value = 0.2116473933682475
$ws.Range("J9").Value = 0.1518081987657302
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.000867
$ws.Range("N9").Value = 0.002601
$ws.Range("O9").Value = 0.7683899556868538
$ws.Range("P9").Value = 0.7683899556868538
$ws.Range("Q9").Value = 0.011600205102
$ws.Range("R9").Value = 0.06960123061199999
$ws.Range("S9").Value = 0.1626277312114658
$ws.Range("T9").Value = 0.1166478951225005

# Row 10
$ws.Range("A10").Value = "Resolving-Mac"
$ws.Range("B10").Value = "Lin7c"
$ws.Range("C10").Value = "Htr2c"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 17.71143666666667
$ws.Range("H10").Value = 53.13431000000001
$ws.Range("I10").Value = 0.2801690413307141
$ws.Range("J10").Value = 0.3014350200878827
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.0002613333333333333
$ws.Range("N10").Value = 0.000784
$ws.Range("O10").Value = 0.2316100443131462
$ws.Range("P10").Value = 0.2316100443131462
$ws.Range("Q10").Value = 0.004628588782222222
$ws.Range("R10").Value = 0.04165729904000001
$ws.Range("S10").Value = 0.06488996407777838
$ws.Range("T10").Value = 0.06981537836008864

# Row 11
$ws.Range("A11").Value = "Resolving-Mac"
$ws.Range("B11").Value = "Lin7c"
$ws.Range("C11").Value = "Htr2c"
$ws.Range("D11").Value = "Inflammatory-Mac"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 17.71143666666667
$ws.Range("H11").Value = 53.13431000000001
$ws.Range("I11").Value = 0.2801690413307141
$ws.Range("J11").Value = 0.3014350200878827
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0.3333333333333333
$ws.Range("M11").Value = 0.000867
$ws.Range("N11").Value = 0.002601
$ws.Range("O11").Value = 0.7683899556868538
$ws.Range("P11").Value = 0.7683899556868538
$ws.Range("Q11").Value = 0.01535581559
$ws.Range("R11").Value = 0.13820234031
$ws.Range("S11").Value = 0.2152790772529357
$ws.Range("T11").Value = 0.2316196417277941
